# Fill in the candidate's answers on the "Mandatory Information Sheet"
# (cells are written in the same order the original author entered them,
# so new shared-string entries land at the same indices as the target file)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = "4 years 4months"
$ws.Range("C3").Value  = "contant learning"
$ws.Range("C4").Value  = "Backend Developer"
$ws.Range("C5").Value  = "Love to tackle real world problems using Data Structures , Algorithms and Computer Science Technologies"
$ws.Range("C6").Value  = "4years 4months"
$ws.Range("C7").Value  = "2months(30-45days)"
$ws.Range("C8").Value  = 7736690165
$ws.Range("C9").Value  = "No"
$ws.Range("C10").Value = "No"
$ws.Range("C11").Value = "B-tech"
$ws.Range("C12").Value = "Single"
$ws.Range("C13").Value = "14lakhs"
$ws.Range("C15").Value = "2.5 lakhs(variable)"
$ws.Range("C16").Value = "50k"
$ws.Range("C26").Value = "Bangalore"
$ws.Range("C18").Value = "    26 lakhs"
$ws.Range("C22").Value = "1)Led client project team of 6 backend developer and delivered on time `n2) Involved in different features development using apis, cloud, db like automating pipeline using api, dockerfile enrichments, onboarding large set of data, `ntop-botoom data mapping and vice-versa, registrations and invoice calculation etc."
$ws.Range("C23").Value = "1) Developed and designed single-click platform to deploy application on different cloud from source code location`n2) Developed and delivered Security management application portal for client `n3) Developing service that Automates purchasing of goods and Services for retails shop without human intervention"
$ws.Range("C25").Value = "Manager"
$ws.Range("C21").Value = "tried to use server resources and written efficient code"
$ws.Range("C24").Value = 3
